# ADD results from server
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # 2025
$ws2 = $wb.Worksheets.Item(2)   # 2030
$ws3 = $wb.Worksheets.Item(3)   # 2035
$ws4 = $wb.Worksheets.Item(4)   # 2040
$ws5 = $wb.Worksheets.Item(5)   # 2045
$ws6 = $wb.Worksheets.Item(6)   # 2050

# Sheet 2025
$ws1.Range("A2").Value = 0
$ws1.Range("B2").Value = 0.07994270341784696
$ws1.Range("E2").Value = 0.1972349740179767
$ws1.Range("I2").Value = 0.8532784
$ws1.Range("L2").Value = 0
$ws1.Range("M2").Value = 0
$ws1.Range("N2").Value = 3.461393823564564
$ws1.Range("O2").Value = 2.471387623027189

# Sheet 2030
$ws2.Range("A2").Value = 0.1909954744387832
$ws2.Range("B2").Value = 0.140464496582153
$ws2.Range("E2").Value = 0.5790800236167917
$ws2.Range("I2").Value = 1.674052599999999
$ws2.Range("L2").Value = 0
$ws2.Range("M2").Value = 0
$ws2.Range("N2").Value = 5.560887517143134
$ws2.Range("O2").Value = 6.19649601055773

# Sheet 2035
$ws3.Range("A2").Value = 0.2539772255612167
$ws3.Range("B2").Value = 0.2296106490751534
$ws3.Range("E2").Value = 0.9196628171525797
$ws3.Range("I2").Value = 2.620466577457396
$ws3.Range("L2").Value = 0
$ws3.Range("M2").Value = 0
$ws3.Range("N2").Value = 23.86138172820945
$ws3.Range("O2").Value = 15.38923335888671

# Sheet 2040
$ws4.Range("O2").Value = 1.863440681153939

# Sheet 2045
$ws5.Range("A2").Value = 0.1161783067483362
$ws5.Range("O2").Value = 1.634994899896032

# Sheet 2050 - unchanged
